$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update prices in D31:D33
$ws.Range("D31").Value = 3985.166
$ws.Range("D32").Value = 1646.582
$ws.Range("D33").Value = 3379.108
